# Add season-record columns (Wins, Losses, Ties) to the player table.
# The sheet currently spans A1:AC53 (player/team stats). We append three
# new columns (AD, AE, AF) holding the team's season record, repeated on
# every player row, with a header row styled like the existing headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column titles ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold font, thin box border,
# centered horizontally, top-aligned vertically) used by columns A1:AC1.
$hdr = $ws.Range("AD1:AF1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.Item(1).LineStyle = 1 # xlContinuous (left)
$hdr.Borders.Item(2).LineStyle = 1 # right
$hdr.Borders.Item(3).LineStyle = 1 # top
$hdr.Borders.Item(4).LineStyle = 1 # bottom

# --- Data rows (2-53): season record values, same for every player row ---
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 72  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 89  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
